$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from the table cell paragraph
#    ("Die Hunde w -> Hunde (die, w)"). Find.Execute repositions the
#    range it is called on to the matched text, so re-use that same
#    range object for the replacement (InsertXML replaces exactly the
#    range it is invoked on, dropping the bookmark that previously sat
#    at the paragraph mark just past the match).
# ---------------------------------------------------------------------
$cellRange = $d.Content
$found = $cellRange.Find.Execute("Die Hunde w -> Hunde (die, w)")
if ($found) {
    $cellXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00450667" w:rsidRDefault="00450667" w:rsidP="005C73D8"><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Die Hunde w -&gt; Hunde (die, w)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $cellRange.InsertXML($cellXml)
}

# ---------------------------------------------------------------------
# 2) Append the new "Begräbnis" discussion paragraphs at the end of the
#    document body (before the final section break), moving the
#    "_GoBack" bookmark into the new last paragraph.
# ---------------------------------------------------------------------
$endRange = $d.Range($d.Content.End, $d.Content.End)
$newParas = '<w:p/>' `
    + '<w:p><w:r><w:t>III_17_0003 Begr' + [char]0x00E4 + 'bnis</w:t></w:r></w:p>' `
    + '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Vergr' + [char]0x00E4 + 'bnis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> als </w:t></w:r><w:r><w:t>zweites</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Wort</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t xml:space="preserve">Was wenn es eine </w:t></w:r><w:r><w:t>Anmerkung</w:t></w:r><w:r><w:t xml:space="preserve"> des </w:t></w:r><w:r><w:t>Informanten</w:t></w:r><w:r><w:t xml:space="preserve"> ist und ich unsicher bin?</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>[&lt;(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vergr' + [char]0x00E4 + 'bnis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)&gt;]</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t xml:space="preserve">Generell </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>was</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> wenn mehrere Klammern zusammenkommen</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t xml:space="preserve">Und was ist wenn etwas vom </w:t></w:r><w:r><w:t>Informanten</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> in () steht</w:t></w:r></w:p>'

$endXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParas + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($endXml)
